# Applies the "output generated at 456a3b4" data refresh to the
# 上海-漫展信息 workbook: updated "want to go" counts (column F) on
# several rows across all four sheets, plus a full row-27 update on
# sheet "展览" (title, availability price, link and cover image).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Sheet 1: 展览 (Exhibitions)
# ---------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("展览")

$ws1.Range("F4").Value  = 371
$ws1.Range("F5").Value  = 8161
$ws1.Range("F7").Value  = 74
$ws1.Range("F8").Value  = 2135
$ws1.Range("F10").Value = 1102
$ws1.Range("F13").Value = 8
$ws1.Range("F14").Value = 1170
$ws1.Range("F16").Value = 21
$ws1.Range("F17").Value = 744
$ws1.Range("F19").Value = 535
$ws1.Range("F20").Value = 64
$ws1.Range("F23").Value = 6856
$ws1.Range("F25").Value = 53958
$ws1.Range("F26").Value = 4193

# Row 27 - event rebranded / relisted for sale
$ws1.Range("C27").Value = "上海·次元空港·千年节·同人动漫游戏嘉年华"
$ws1.Range("F27").Value = 1
$ws1.Range("G27").Value = 68.2
$ws1.Range("H27").Value = "https://show.bilibili.com/platform/detail.html?id=92403"
$ws1.Range("I27").Value = "//i0.hdslb.com/bfs/openplatform/202409/yvHkiDMK1726117265087.jpeg"

$ws1.Range("F28").Value = 1015
$ws1.Range("F29").Value = 801
$ws1.Range("F30").Value = 382
$ws1.Range("F35").Value = 2035
$ws1.Range("F39").Value = 1079
$ws1.Range("F40").Value = 469
$ws1.Range("F42").Value = 169
$ws1.Range("F45").Value = 122
$ws1.Range("F47").Value = 116

# ---------------------------------------------------------------
# Sheet 2: 演出 (Performances)
# ---------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("演出")

$ws2.Range("F9").Value  = 137
$ws2.Range("F12").Value = 41
$ws2.Range("F13").Value = 96
$ws2.Range("F15").Value = 31
$ws2.Range("F17").Value = 7376
$ws2.Range("F18").Value = 94
$ws2.Range("F26").Value = 16
$ws2.Range("F28").Value = 112

# ---------------------------------------------------------------
# Sheet 3: 本地生活 (Local life)
# ---------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("本地生活")

$ws3.Range("F4").Value  = 2246
$ws3.Range("F5").Value  = 1504
$ws3.Range("F7").Value  = 642
$ws3.Range("F8").Value  = 2331
$ws3.Range("F9").Value  = 9286
$ws3.Range("F10").Value = 1577
$ws3.Range("F15").Value = 111

# ---------------------------------------------------------------
# Sheet 4: 全部类型 (All types)
# ---------------------------------------------------------------
$ws4 = $wb.Worksheets.Item("全部类型")

$ws4.Range("F4").Value  = 2246
$ws4.Range("F5").Value  = 371
$ws4.Range("F6").Value  = 642
$ws4.Range("F7").Value  = 1577
$ws4.Range("F10").Value = 74
$ws4.Range("F13").Value = 8
$ws4.Range("F14").Value = 1170
$ws4.Range("F15").Value = 111
$ws4.Range("F16").Value = 21
$ws4.Range("F17").Value = 744
$ws4.Range("F18").Value = 64
$ws4.Range("F20").Value = 6856
$ws4.Range("F22").Value = 53958
$ws4.Range("F23").Value = 137
$ws4.Range("F24").Value = 137
$ws4.Range("F27").Value = 4193
$ws4.Range("F28").Value = 1015
$ws4.Range("F29").Value = 382
$ws4.Range("F33").Value = 96
$ws4.Range("F36").Value = 31
$ws4.Range("F39").Value = 94
$ws4.Range("F42").Value = 169
$ws4.Range("F44").Value = 122
$ws4.Range("F46").Value = 116
